$d = $word.ActiveDocument

# --- 1) Clear leftover placeholder "John" (default student first name) ---
# Appears twice as filler text next to "Nom :" labels; blank it out both times.
$d.Content.Find.Execute("John", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# --- 2) Clear leftover placeholder "MaitreDeStageNom" (default internship supervisor name) ---
# Appears twice: once as plain text after "Nom :  ", and once as the *result*
# of a " MERGEFIELD PRENOM_ENCADRANT " field. A plain Find/Replace across the
# document body only reaches the non-field occurrence, so handle that first...
$d.Content.Find.Execute("MaitreDeStageNom", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ...then locate the field whose displayed result still contains the
# placeholder and strip the word from its result range directly, leaving the
# leading space intact (" MaitreDeStageNom" -> " ").
$fields = $d.Fields
for ($i = 1; $i -le $fields.Count; $i++) {
    $f = $fields.Item($i)
    if ($f.Result.Text -like "*MaitreDeStageNom*") {
        $resultRange = $f.Result
        $wordStart = $resultRange.Start + ($resultRange.Text.Length - "MaitreDeStageNom".Length)
        $wordRange = $d.Range($wordStart, $resultRange.End)
        $wordRange.Delete()
    }
}

# --- 3) Replace leftover placeholder "TuteurNom" (default academic tutor name) with "null" ---
# Appears twice, both times next to an already-blank "null" first-name merge
# field, so the result reads "null null" in both spots.
$d.Content.Find.Execute("TuteurNom", $true, $false, $false, $false, $false, $true, 1, $false, "null", 2)
